{"js": "// Vietnamese translation pass for\n// \"Email 5-3 [TEMPLATE] Partner email \u2013 document verification failed\"\n//\n// Each change below replaces an English run's text with its Vietnamese\n// counterpart, matching the unified diff 1:1 (run-scoped text swaps only,\n// no structural changes to paragraphs/runs/formatting).\n\nconst body = context.document.body;\n\n// Replace the first occurrence of `find` within `searchScope` (a Body or\n// Range) with `replaceWith`, preserving the run's own formatting.\nasync function replaceOnce(searchScope, find, replaceWith) {\n  const results = searchScope.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n  results.items[0].insertText(replaceWith, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Locate the two paragraphs whose \" or \" runs are otherwise ambiguous at the\n// body level (there are two identical \" or \" runs in the document), so we\n// resolve them once by their distinctive paragraph text and reuse the stable\n// paragraph reference for every edit inside that paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet liveChatPara = null;\nlet countryManagerPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"If you have any questions, please contact us via\") === 0) {\n    liveChatPara = paragraphs.items[i];\n  } else if (t.indexOf(\"If you have any questions, please contact your country manager,\") === 0) {\n    countryManagerPara = paragraphs.items[i];\n  }\n}\nif (!liveChatPara || !countryManagerPara) {\n  throw new Error(\"Could not locate target paragraphs\");\n}\nconst liveChatRange = liveChatPara.getRange();\nconst countryManagerRange = countryManagerPara.getRange();\n\n// 1) Heading\nawait replaceOnce(body, \"Uh oh! Your documents couldn\\u2019t be verified\", \"R\u1ea5t ti\u1ebfc! Ch\u00fang t\u00f4i kh\u00f4ng th\u1ec3 x\u00e1c minh th\u00f4ng tin c\u1ee7a b\u1ea1n\");\n\n// 2) Greeting\nawait replaceOnce(body, \"Hi \", \"Xin ch\u00e0o \");\n\n// 3) Intro paragraph\nawait replaceOnce(\n  body,\n  \"We regret to inform you that your documents have failed our verification process as we found the following issues with them: \",\n  \"Ch\u00fang t\u00f4i r\u1ea5t ti\u1ebfc ph\u1ea3i th\u00f4ng b\u00e1o \u0111\u1ebfn b\u1ea1n r\u1eb1ng ch\u00fang t\u00f4i kh\u00f4ng th\u1ec3 x\u00e1c minh th\u00f4ng tin b\u1ea1n g\u1eedi v\u00ec m\u1ed9t s\u1ed1 v\u1ea5n \u0111\u1ec1 sau \u0111\u00e2y: \"\n);\n\n// 4) Bullet item \u2014 bold label\nawait replaceOnce(body, \"A copy of your vaccination certificate\", \"B\u1ea3n sao gi\u1ea5y ch\u1ee9ng nh\u1eadn ti\u00eam ch\u1ee7ng c\u1ee7a b\u1ea1n\");\n\n// 5) Bullet item \u2014 explanation\nawait replaceOnce(body, \": Document is unclear\", \": Gi\u1ea5y t\u1edd kh\u00f4ng r\u00f5 r\u00e0ng\");\n\n// 6) Resubmit instructions\nawait replaceOnce(body, \"Please resubmit the documents above by \", \"Xin vui l\u00f2ng g\u1eedi l\u1ea1i c\u00e1c gi\u1ea5y t\u1edd tr\u00ean tr\u01b0\u1edbc ng\u00e0y \");\n\n// 7) Tail of resubmit sentence\nawait replaceOnce(body, \" so we can proceed with the necessary arrangements.\", \" \u0111\u1ec3 ch\u00fang t\u00f4i c\u00f3 th\u1ec3 ti\u1ebfp t\u1ee5c c\u00e1c b\u01b0\u1edbc s\u1eafp x\u1ebfp c\u1ea7n thi\u1ebft.\");\n\n// 8) \"contact us via\" paragraph intro\nawait replaceOnce(liveChatRange, \"If you have any questions, please contact us via \", \"N\u1ebfu b\u1ea1n c\u1ea7n h\u1ed7 tr\u1ee3, h\u00e3y li\u00ean h\u1ec7 v\u1edbi ch\u00fang t\u00f4i qua \");\n\n// 9) \" or \" between live chat / WhatsApp links\nawait replaceOnce(liveChatRange, \" or \", \" ho\u1eb7c \");\n\n// 10) \"country manager\" paragraph intro\nawait replaceOnce(countryManagerRange, \"If you have any questions, please contact your country manager, \", \"N\u1ebfu b\u1ea1n c\u00f3 b\u1ea5t k\u1ef3 th\u1eafc m\u1eafc n\u00e0o, vui l\u00f2ng li\u00ean h\u1ec7 v\u1edbi gi\u00e1m \u0111\u1ed1c ph\u1ee5 tr\u00e1ch qu\u1ed1c gia c\u1ee7a b\u1ea1n \");\n\n// 11) \", at \" -> \", qua email \"\nawait replaceOnce(countryManagerRange, \", at \", \", qua email \");\n\n// 12) \" or \" before WhatsApp number -> \" ho\u1eb7c s\u1ed1 \"\nawait replaceOnce(countryManagerRange, \" or \", \" ho\u1eb7c s\u1ed1 \");\n", "ps1": "# Vietnamese translation pass for\n# \"Email 5-3 [TEMPLATE] Partner email - document verification failed\"\n#\n# Each change below replaces an English run's text with its Vietnamese\n# counterpart, matching the unified diff 1:1 (text-only swaps, no structural\n# changes to paragraphs/runs/formatting).\n#\n# NOTE: this interpreter only binds user-defined function parameters\n# positionally, so Replace-InRange is always invoked with positional args\n# (Range, FindText, ReplaceWith) rather than -Name value pairs.\n\n$d = $word.ActiveDocument\n\nfunction Replace-InRange {\n    param($Range, [string]$FindText, [string]$ReplaceWith)\n    $find = $Range.Find\n    $find.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceWith\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdFindContinue = 1 (Wrap), wdReplaceOne = 1 (Replace only the first hit)\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n}\n\n# 1) Heading\nReplace-InRange $d.Content \"Uh oh! Your documents couldn't be verified\" \"R\u1ea5t ti\u1ebfc! Ch\u00fang t\u00f4i kh\u00f4ng th\u1ec3 x\u00e1c minh th\u00f4ng tin c\u1ee7a b\u1ea1n\"\n\n# 2) Greeting\nReplace-InRange $d.Content \"Hi \" \"Xin ch\u00e0o \"\n\n# 3) Intro paragraph\nReplace-InRange $d.Content \"We regret to inform you that your documents have failed our verification process as we found the following issues with them: \" \"Ch\u00fang t\u00f4i r\u1ea5t ti\u1ebfc ph\u1ea3i th\u00f4ng b\u00e1o \u0111\u1ebfn b\u1ea1n r\u1eb1ng ch\u00fang t\u00f4i kh\u00f4ng th\u1ec3 x\u00e1c minh th\u00f4ng tin b\u1ea1n g\u1eedi v\u00ec m\u1ed9t s\u1ed1 v\u1ea5n \u0111\u1ec1 sau \u0111\u00e2y: \"\n\n# 4) Bullet item - bold label\nReplace-InRange $d.Content \"A copy of your vaccination certificate\" \"B\u1ea3n sao gi\u1ea5y ch\u1ee9ng nh\u1eadn ti\u00eam ch\u1ee7ng c\u1ee7a b\u1ea1n\"\n\n# 5) Bullet item - explanation\nReplace-InRange $d.Content \": Document is unclear\" \": Gi\u1ea5y t\u1edd kh\u00f4ng r\u00f5 r\u00e0ng\"\n\n# 6) Resubmit instructions\nReplace-InRange $d.Content \"Please resubmit the documents above by \" \"Xin vui l\u00f2ng g\u1eedi l\u1ea1i c\u00e1c gi\u1ea5y t\u1edd tr\u00ean tr\u01b0\u1edbc ng\u00e0y \"\n\n# 7) Tail of resubmit sentence\nReplace-InRange $d.Content \" so we can proceed with the necessary arrangements.\" \" \u0111\u1ec3 ch\u00fang t\u00f4i c\u00f3 th\u1ec3 ti\u1ebfp t\u1ee5c c\u00e1c b\u01b0\u1edbc s\u1eafp x\u1ebfp c\u1ea7n thi\u1ebft.\"\n\n# Locate the two paragraphs whose \" or \" / intro runs are otherwise ambiguous\n# at the document level (two identical \" or \" runs exist), and scope those\n# edits to the correct paragraph's Range.\n$liveChatPara = $null\n$countryManagerPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"If you have any questions, please contact us via\")) {\n        $liveChatPara = $p\n    } elseif ($t.StartsWith(\"If you have any questions, please contact your country manager,\")) {\n        $countryManagerPara = $p\n    }\n}\nif ($liveChatPara -eq $null -or $countryManagerPara -eq $null) {\n    throw \"Could not locate target paragraphs\"\n}\n\n# 8) \"contact us via\" paragraph intro\nReplace-InRange $liveChatPara.Range \"If you have any questions, please contact us via \" \"N\u1ebfu b\u1ea1n c\u1ea7n h\u1ed7 tr\u1ee3, h\u00e3y li\u00ean h\u1ec7 v\u1edbi ch\u00fang t\u00f4i qua \"\n\n# 9) \" or \" between live chat / WhatsApp links\nReplace-InRange $liveChatPara.Range \" or \" \" ho\u1eb7c \"\n\n# 10) \"country manager\" paragraph intro\nReplace-InRange $countryManagerPara.Range \"If you have any questions, please contact your country manager, \" \"N\u1ebfu b\u1ea1n c\u00f3 b\u1ea5t k\u1ef3 th\u1eafc m\u1eafc n\u00e0o, vui l\u00f2ng li\u00ean h\u1ec7 v\u1edbi gi\u00e1m \u0111\u1ed1c ph\u1ee5 tr\u00e1ch qu\u1ed1c gia c\u1ee7a b\u1ea1n \"\n\n# 11) \", at \" -> \", qua email \"\nReplace-InRange $countryManagerPara.Range \", at \" \", qua email \"\n\n# 12) \" or \" before WhatsApp number -> \" ho\u1eb7c s\u1ed1 \"\nReplace-InRange $countryManagerPara.Range \" or \" \" ho\u1eb7c s\u1ed1 \"\n"}
